$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 1.86

$ws.Range("F3").Value = 1.92
$ws.Range("K3").Value = 3.95

$ws.Range("J5").Value = 3.95

$ws.Range("G6").Value = 3.4
$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 3.5
$ws.Range("P6").Value = 1.91
$ws.Range("V6").Value = 1.66
$ws.Range("W6").Value = 1.42
$ws.Range("AO6").Value = 23

$ws.Range("F7").Value = 1.42
$ws.Range("G7").Value = 1.49
$ws.Range("K7").Value = 5.4
$ws.Range("P7").Value = 2.22
$ws.Range("Q7").Value = 1.66
$ws.Range("U7").Value = 1.97

$ws.Range("J8").Value = 3.4

$ws.Range("G9").Value = 2.24
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 4.8
$ws.Range("J9").Value = 3.15
$ws.Range("K9").Value = 4.6
$ws.Range("N9").Value = 2.76
$ws.Range("Q9").Value = 1.72
$ws.Range("W9").Value = 1.81
